$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.041.08'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.874.15'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.01%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.81'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.40%  '

$ws.Range('E6').Value = '  +0.28%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5032'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3827'
$ws.Range('D8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08383'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -9.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.111'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.40%  '

$ws.Range('B11').Value = 'Polkadot'
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.231'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.25%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.881.61'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.92%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.49'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.48%  '

$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.199'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.16%  '

$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.006'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.41%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001096'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.07%  '

$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.01'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.53%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06667'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.17%  '

$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.97'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.89%  '

$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.003'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.23%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.065'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.21%  '

$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '28.089.87'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.18%  '

$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.15'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.26%  '

$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.272'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.19%  '

$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.574'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.21%  '

$ws.Range('B26').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C26').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.098.31'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.70%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.62'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.03%  '

$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '156.42'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.11%  '

$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '126.54'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.28%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1047'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.92%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.047'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.57%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.615'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.16%  '

$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.611'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('B34').Value = 'FraxShare'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.683'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.99%  '

$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02446'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.84%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06525'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.26%  '

$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2167'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.77%  '

$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.218'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.02%  '

$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6527'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.04%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.244'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -8.15%  '

$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.31'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.99%  '

$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.883'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.66%  '

$ws.Range('B43').Value = 'Decentraland'
$ws.Range('C43').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6214'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.31%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.03'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.65%  '

$ws.Range('B45').Value = 'WEMIXTOKEN'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.300'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.38%  '

$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.681'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.05%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.009'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.220'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.83%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '120.77'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.05%  '

$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '80.18'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.49%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06894'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.35%  '
